# Update the "Date Placeholder 1" shape on every slide from 8/6/2019 to
# 8/9/2019, and make sure its position/size is explicitly pinned to the
# standard placeholder frame (off 628650,6356351 ext 2057400,365125 EMU,
# i.e. 49.5,500.5001,162,28.75 points) instead of inheriting from the
# layout/master.

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.Name -eq "Date Placeholder 1") {
            $shp.Left = 49.5
            $shp.Top = 500.5001
            $shp.Width = 162
            $shp.Height = 28.75
            if ($shp.TextFrame.TextRange.Text -eq "8/6/2019") {
                $shp.TextFrame.TextRange.Text = "8/9/2019"
            }
        }
    }
}
